# add save column in s_vals sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header - copy format from G1 (same header style) then set the text
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values per row (2-14)
$saveValues = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 1
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
